# This script applies the "environmental.docx" content edit described by the
# commit "added additional changes to contents for website".
#
# Word COM (and this headless iron_native shim) does not expose a way to
# programmatically split an existing run into several runs or to insert
# brand-new <w:proofErr/> markers through the normal Range/Find API, so the
# most faithful and robust way to reproduce the target OOXML is to replace
# the two affected paragraphs' Range content with literal WordprocessingML
# via Range.InsertXML(), matching exactly the run/proofErr/pPr structure
# that Word itself produced for this edit.

$d = $word.ActiveDocument

# --- Paragraph 3 ("One of the leading problem here in Baguio City...") ---
# The only change here is splitting "...tragic event happened in " into
# three runs so the word "that" is inserted as its own run:
#   "...tragic event " + "that " + "happened in "
$para3Xml = @'
<w:p>
<w:r>
<w:t xml:space="preserve">One of the leading problem here in Baguio City is Garbage Collection. </w:t>
</w:r>
<w:r>
<w:t xml:space="preserve">Do you still remember the tragic event </w:t>
</w:r>
<w:r>
<w:t xml:space="preserve">that </w:t>
</w:r>
<w:r>
<w:t xml:space="preserve">happened in </w:t>
</w:r>
<w:proofErr w:type="spellStart"/>
<w:r>
<w:t>Barangay</w:t>
</w:r>
<w:proofErr w:type="spellEnd"/>
<w:r>
<w:t xml:space="preserve"> </w:t>
</w:r>
<w:proofErr w:type="spellStart"/>
<w:r>
<w:t>Asin</w:t>
</w:r>
<w:proofErr w:type="spellEnd"/>
<w:r>
<w:t xml:space="preserve">? Maybe few would recognize that story today. Baguio residents remember </w:t>
</w:r>
<w:proofErr w:type="spellStart"/>
<w:r>
<w:t>Asin</w:t>
</w:r>
<w:proofErr w:type="spellEnd"/>
<w:r>
<w:t xml:space="preserve"> Road as the neighbouring residential district on which the city old dump collapsed in year 2011.</w:t>
</w:r>
</w:p>
'@
$r3 = $d.Paragraphs(3).Range
$null = $r3.InsertXML($para3Xml)

# --- Paragraph 4 ("Barangay Asin Officials are aware...") and everything  ---
# --- that follows it (new trailing paragraphs) ------------------------------
# This paragraph gains a new <w:pPr> (single spacing, no space-after, a
# first-line indent), several runs are re-split/re-worded ("fined" ->
# "fine" + "d. ...", "barangay" -> "Barangay" with spell-check markers,
# "clogged in the canals" -> "clogging", "Clean and Green's" -> "Clean and
# Green", an extra "the " before "bayanihan"), and the final sentence
# ("Having trashcans...") is replaced by several new sentences. Three new
# paragraphs (one blank, one about "kagawads", one about "within their...")
# are appended after it, followed by one more trailing blank paragraph.
$restXml = @'
<w:p>
<w:pPr>
<w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
<w:ind w:firstLine="720"/>
</w:pPr>
<w:proofErr w:type="spellStart"/>
<w:r>
<w:t>Barangay</w:t>
</w:r>
<w:proofErr w:type="spellEnd"/>
<w:r>
<w:t xml:space="preserve"> </w:t>
</w:r>
<w:proofErr w:type="spellStart"/>
<w:r>
<w:t>Asin</w:t>
</w:r>
<w:proofErr w:type="spellEnd"/>
<w:r>
<w:t xml:space="preserve"> Officials are aware of the problem and are making an effort in helping to clean up the mess. The community have started segregating household waste. It is really strictly implemented that if you don’t abide by the rules you’ll be fine</w:t>
</w:r>
<w:r>
<w:t xml:space="preserve">d. One of the committee in the </w:t>
</w:r>
<w:proofErr w:type="spellStart"/>
<w:r>
<w:t>B</w:t>
</w:r>
<w:r>
<w:t>arangay</w:t>
</w:r>
<w:proofErr w:type="spellEnd"/>
<w:r>
<w:t xml:space="preserve"> is the Clea</w:t>
</w:r>
<w:r>
<w:t>n and Green, the official makes sure</w:t>
</w:r>
<w:r>
<w:t xml:space="preserve"> that every house hold </w:t>
</w:r>
<w:r>
<w:t>are following the rules and that the</w:t>
</w:r>
<w:r>
<w:t xml:space="preserve"> canals are </w:t>
</w:r>
<w:r>
<w:t>garbage-free</w:t>
</w:r>
<w:r>
<w:t xml:space="preserve"> to prev</w:t>
</w:r>
<w:r>
<w:t xml:space="preserve">ent clogging. The </w:t>
</w:r>
<w:proofErr w:type="spellStart"/>
<w:r>
<w:t>B</w:t>
</w:r>
<w:r>
<w:t>arangay</w:t>
</w:r>
<w:proofErr w:type="spellEnd"/>
<w:r>
<w:t xml:space="preserve"> organizes Clean and Green together with </w:t>
</w:r>
<w:r>
<w:t xml:space="preserve">the </w:t>
</w:r>
<w:proofErr w:type="spellStart"/>
<w:r>
<w:t>bayanihan</w:t>
</w:r>
<w:proofErr w:type="spellEnd"/>
<w:r>
<w:t xml:space="preserve"> of the whole community. </w:t>
</w:r>
<w:r>
<w:t xml:space="preserve">Since sanitation is a priority in the </w:t>
</w:r>
<w:proofErr w:type="spellStart"/>
<w:r>
<w:t>Barangay</w:t>
</w:r>
<w:proofErr w:type="spellEnd"/>
<w:r>
<w:t xml:space="preserve">, Street sweepers are the ones responsible for cleaning the different </w:t>
</w:r>
<w:proofErr w:type="spellStart"/>
<w:r>
<w:t>purok</w:t>
</w:r>
<w:proofErr w:type="spellEnd"/>
<w:r>
<w:t xml:space="preserve"> of the </w:t>
</w:r>
<w:proofErr w:type="spellStart"/>
<w:r>
<w:t>Barangay</w:t>
</w:r>
<w:proofErr w:type="spellEnd"/>
<w:r>
<w:t>. They also provide</w:t>
</w:r>
<w:r>
<w:t xml:space="preserve"> trashcans along the road </w:t>
</w:r>
<w:r>
<w:t xml:space="preserve">which </w:t>
</w:r>
<w:r>
<w:t>really helps the community’s cleanliness.</w:t>
</w:r>
<w:bookmarkStart w:id="0" w:name="_GoBack"/>
<w:bookmarkEnd w:id="0"/>
</w:p>
<w:p>
<w:pPr>
<w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
</w:pPr>
</w:p>
<w:p>
<w:pPr>
<w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
</w:pPr>
<w:r>
<w:t xml:space="preserve">According to one of their </w:t>
</w:r>
<w:proofErr w:type="spellStart"/>
<w:r>
<w:t>kagawads</w:t>
</w:r>
<w:proofErr w:type="spellEnd"/>
<w:r>
<w:t xml:space="preserve">, there are also issues </w:t>
</w:r>
</w:p>
<w:p>
<w:pPr>
<w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
</w:pPr>
<w:proofErr w:type="gramStart"/>
<w:r>
<w:t>within</w:t>
</w:r>
<w:proofErr w:type="gramEnd"/>
<w:r>
<w:t xml:space="preserve"> their. </w:t>
</w:r>
<w:proofErr w:type="gramStart"/>
<w:r>
<w:t xml:space="preserve">Issues such as residents having no discipline within themselves in order to protect the </w:t>
</w:r>
<w:proofErr w:type="spellStart"/>
<w:r>
<w:t>Barangay's</w:t>
</w:r>
<w:proofErr w:type="spellEnd"/>
<w:r>
<w:t xml:space="preserve"> environment.</w:t>
</w:r>
<w:proofErr w:type="gramEnd"/>
</w:p>
<w:p>
<w:pPr>
<w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
</w:pPr>
</w:p>
'@
$r4 = $d.Paragraphs(4).Range
$null = $r4.InsertXML($restXml)
